$d = $word.ActiveDocument

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs($index)
    $s = $p.Range.Start
    $e = $p.Range.End
    $d.Range($s, $e).Text = $text
}

function Remove-Para($index) {
    $p = $d.Paragraphs($index)
    $s = $p.Range.Start
    $e = $p.Range.End
    $d.Range($s, $e).Delete()
}

# --- 1. Paragraph "Квалификационная работа студента..." ---
# Remove "нейронной сети " between "модели" and "на программируемых"
$d.Content.Find.Execute(
    "реализации модели нейронной сети на программируемых",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "реализации модели на программируемых", 2)

# --- 2. Paragraph (was) "В рамках данной исследовательской работы были рассмотрены..." ---
Set-ParaText 5 "В рамках данной исследовательской работы было осуществлено создание и конфигурирование нейронной сети, разработанной для решения задачи распознавания лиц на изображениях. Значительное внимание уделено подбору и обучению оптимальной модели. Проведён выбор её архитектуры, определение параметров и настройка процесса обучения. Результаты исследования позволили достичь высокой точности распознавания лиц на изображениях."

# --- 3. Paragraph (was) "Во-вторых, значительное внимание уделено..." ---
Set-ParaText 6 "Также было проведено тестирование разработанной системы с целью оценки ее производительности и эффективности в реальных условиях."

# --- 4. Paragraph (was) "В третьей части работы была реализована..." ---
Set-ParaText 7 "Объем расчетно-пояснительной записки – __ листов формата А4, графическая часть – __ листов формата А1."

# Word leaves a "_GoBack" bookmark at the last edited spot - place it right
# after "А1" (before the closing period), matching where the author's
# cursor last made changes.
$pVol = $d.Paragraphs(7)
$volStart = $pVol.Range.Start
$volText = $pVol.Range.Text
$marker = "А1."
$markerIdx = $volText.IndexOf($marker)
if ($markerIdx -ge 0) {
    $bmPos = $volStart + $markerIdx + 2
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- 5. Delete paragraph (was) "Исследовательская работа представляет собой комплексный подход..." ---
Remove-Para 8

# --- 6. Delete paragraph (was) "Объем расчетно-пояснительной записки – __ листов формата А4, графическая часть – __ листов формата А4." ---
Remove-Para 8

# --- 7. Paragraph (was) "К достоинствам данной работы следует отнести..." ---
Set-ParaText 8 "К достоинствам данной работы следует отнести ее практическую значимость: разработка нейронной сети для распознавания лиц на изображениях является актуальной и востребованной задачей в современных информационных системах и приложениях."

# --- 8. Paragraph (was) "Однако, в работе также имеются некоторые замечания..." ---
Set-ParaText 9 "К замечаниям стоит отнести пожелание улучшить обоснование выбора конкретной модели нейронной сети и подбор используемых параметров ее обучения."

# --- 9. Paragraph (was) "В целом, выполненная работа соответствует..." ---
$d.Content.Find.Execute(
    "В целом, выполненная работа соответствует",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Выполненная работа соответствует", 2)

Write-Host "DONE"
